# Reorder the "Recorded By" names/emails in column G for the specific rows
# that changed (author/recorder list order swapped), as described by the
# diff. Only the text content of these cells changes - no style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    3   = "System, dnasr281@gmail.com"
    5   = "backup@backdoor.com, System"
    6   = "System, dnasr281@gmail.com"
    7   = "admin@admin.com, System"
    8   = "backup@backdoor.com, System"
    28  = "backup@backdoor.com, System, system"
    29  = "System, dnasr281@gmail.com"
    31  = "backup@backdoor.com, System"
    32  = "System, dnasr281@gmail.com"
    33  = "admin@admin.com, System"
    34  = "backup@backdoor.com, System"
    54  = "backup@backdoor.com, System, system"
    55  = "System, dnasr281@gmail.com"
    57  = "backup@backdoor.com, System"
    58  = "System, dnasr281@gmail.com"
    59  = "admin@admin.com, System"
    60  = "backup@backdoor.com, System"
    80  = "backup@backdoor.com, System"
    81  = "backup@backdoor.com, System"
    82  = "backup@backdoor.com, System"
    87  = "admin@admin.com, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    113 = "admin@admin.com, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    139 = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
